$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial value of 45189 (2023-09-20)
# for every data row (rows 2-359). Update it to 45190 (2023-09-21).
$range = $ws.Range("C2:C359")
$range.Value = 45190
